$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.599.19'
$ws.Range("E2").Value = '  +1.34%  '

$ws.Range("D3").Value = '2.583.92'
$ws.Range("E3").Value = '  +9.92%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.15%  '

$ws.Range("E7").Value = '  +5.73%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +12.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.63'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +12.21%  '

$ws.Range("E11").Value = '  +5.29%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.15'
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").Value = '2.979.28'
$ws.Range("E13").Value = '  +9.90%  '

$ws.Range("E14").Value = '  +1.85%  '

$ws.Range("D15").Value = '2.584.61'
$ws.Range("E15").Value = '  +8.61%  '

$ws.Range("E16").Value = '  +12.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.92'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +9.55%  '

$ws.Range("D18").Value = '46.760.96'
$ws.Range("E18").Value = '  +1.92%  '

$ws.Range("E19").Value = '  +5.59%  '

$ws.Range("E20").Value = '  +4.37%  '

$ws.Range("E21").Value = '  +11.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.78'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '254.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.97%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.03%  '

$ws.Range("E25").Value = '  +14.56%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.06'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +34.55%  '

$ws.Range("E27").Value = '  +0.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.50'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.18%  '

$ws.Range("E29").Value = '  +4.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '39.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.76'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +11.33%  '

$ws.Range("E33").Value = '  +24.81%  '

$ws.Range("E34").Value = '  +5.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0830'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.38%  '

$ws.Range("E36").Value = '  +2.71%  '

$ws.Range("E37").Value = '  +4.31%  '

$ws.Range("E38").Value = '  +5.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.19'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '15.72'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.54%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +12.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0324'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.82%  '

$ws.Range("D43").Value = '2.018.76'
$ws.Range("E43").Value = '  +7.52%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.43'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +30.56%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.26'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.80'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.09'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '108.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +11.59%  '

$ws.Range("E50").Value = '  +7.86%  '

$ws.Range("D51").Value = '2.837.31'
$ws.Range("E51").Value = '  +9.88%  '
